$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
try {
$ws.Rows.Item(18).AutoFit()
Write-Output "autofit ok"
} catch {
Write-Output "autofit ERR: $_"
}
